$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (AJAX): C5, D5, E5 -> -441003.5953130126
$ws.Range("C5").Value = -441003.5953130126
$ws.Range("D5").Value = -441003.5953130126
$ws.Range("E5").Value = -441003.5953130126

# Row 17 (EGYPT): C17, D17 -> -4134175.702280757 ; E17 -> -4134175.702280753
$ws.Range("C17").Value = -4134175.702280757
$ws.Range("D17").Value = -4134175.702280757
$ws.Range("E17").Value = -4134175.702280753

# Row 22 (IBM1): D22 -> 287.1054198927933
$ws.Range("D22").Value = 287.1054198927933

# Row 40 (PDI): C40 -> -294070 ; D40 -> -294070 ; E40 -> -294069.9999999999
$ws.Range("C40").Value = -294070
$ws.Range("D40").Value = -294070
$ws.Range("E40").Value = -294069.9999999999

# Row 54 (TABORA): D54 -> -8471.957300271906
$ws.Range("D54").Value = -8471.957300271906
